# Weekly update: add this week's two new price rows (Primera / Segunda)
# for "Vega Monumental Concepción - Zanahoria" just above the prior-week
# rows, shifting all existing data rows from 349:373 down to 351:375.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows above the current row 349 (first data row of
# the most-recent-week block), pushing everything else down by two rows.
$ws.Rows("349:350").Insert()

# --- New row 349 : Calidad = Primera -------------------------------------
$ws.Range("A349").Value = 11
$ws.Range("B349").Value = "Vega Monumental Concepción"
$ws.Range("C349").Value = "Bíobío"
$ws.Range("D349").Value = 45021
$ws.Range("E349").Value = 8
$ws.Range("F349").Value = 100114013
$ws.Range("G349").Value = "Zanahoria"
$ws.Range("H349").Value = "Sin especificar"
$ws.Range("I349").Value = "Primera"
$ws.Range("J349").Value = 600
$ws.Range("K349").Value = 4500
$ws.Range("L349").Value = 5000
$ws.Range("M349").Value = 4750
$ws.Range("N349").Value = "$/saco 20 kilos"
$ws.Range("O349").Value = "Región de Ñuble"
$ws.Range("P349").Value = 238
$ws.Range("Q349").Value = 20
$ws.Range("R349").Value = "Hortaliza"

# --- New row 350 : Calidad = Segunda -------------------------------------
$ws.Range("A350").Value = 11
$ws.Range("B350").Value = "Vega Monumental Concepción"
$ws.Range("C350").Value = "Bíobío"
$ws.Range("D350").Value = 45021
$ws.Range("E350").Value = 8
$ws.Range("F350").Value = 100114013
$ws.Range("G350").Value = "Zanahoria"
$ws.Range("H350").Value = "Sin especificar"
$ws.Range("I350").Value = "Segunda"
$ws.Range("J350").Value = 300
$ws.Range("K350").Value = 4000
$ws.Range("L350").Value = 4000
$ws.Range("M350").Value = 4000
$ws.Range("N350").Value = "$/saco 20 kilos"
$ws.Range("O350").Value = "Región de Ñuble"
$ws.Range("P350").Value = 200
$ws.Range("Q350").Value = 20
$ws.Range("R350").Value = "Hortaliza"
